$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# The newest quarter's row (02/13/2025, Q4 24, 12/24) is being removed from
# the earnings history table; delete it and shift everything else up.
$pic = $ws.Shapes.Item(1)
$picTop = $pic.Top
$picHeight = $pic.Height

$ws.Rows.Item(2).Delete()

# The floating picture doesn't automatically track the row shift, so nudge
# it up by one default row height (15pt) to keep it anchored the same way
# relative to the remaining rows.
$pic.Top = $picTop - 15
$pic.Height = $picHeight

# Reflect the post-edit selection.
$ws.Activate()
$ws.Range("C8").Select()
